# Daily auto push: a new day of scraped data (2026/01/26 -> 2026/01/27)
# gets inserted into the middle of the log, just before the 2026/12/29
# block, shifting everything below it down by two rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 730; the old rows 730..771 become 732..773.
$ws.Rows.Item(730).Resize(2).Insert()

# New row 730: 2026/01/26 (Mon), time 22, rank 174.
$c = $ws.Cells.Item(730, 1)
$c.NumberFormat = "@"
$c.Value = "2026/01/26"
$c.ClearFormats()
$ws.Cells.Item(730, 2).Value = "月"
$ws.Cells.Item(730, 3).Value = 22
$ws.Cells.Item(730, 4).Value = 174

# New row 731: 2026/01/27 (Tue), time 1, rank 187.
$c = $ws.Cells.Item(731, 1)
$c.NumberFormat = "@"
$c.Value = "2026/01/27"
$c.ClearFormats()
$ws.Cells.Item(731, 2).Value = "火"
$ws.Cells.Item(731, 3).Value = 1
$ws.Cells.Item(731, 4).Value = 187
